# Refresh the crypto price/volume snapshot (GitHub Actions daily update).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price/volume columns hold numeric-looking text (e.g. "39.530.34", "0.636")
# that must stay text, not be coerced into Excel numbers. Force the cell to
# Text format while writing the value, then restore the default "Normal"
# style so no stray number-format styling is left behind on the cell.
function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

# Row 2
Set-TextValue $ws.Range('D2') '39.530.34'
Set-TextValue $ws.Range('E2') '  +1.99%  '

# Row 3
Set-TextValue $ws.Range('D3') '2.168.29'
Set-TextValue $ws.Range('E3') '  +3.18%  '

# Row 4
Set-TextValue $ws.Range('E4') '  +0.01%  '

# Row 5
Set-TextValue $ws.Range('D5') '228.94'
Set-TextValue $ws.Range('E5') '  +0.58%  '

# Row 6
Set-TextValue $ws.Range('D6') '0.636'
Set-TextValue $ws.Range('E6') '  +3.12%  '

# Row 7
Set-TextValue $ws.Range('D7') '63.72'
Set-TextValue $ws.Range('E7') '  +2.39%  '

# Row 8
Set-TextValue $ws.Range('E8') '  +0.02%  '

# Row 9
Set-TextValue $ws.Range('E9') '  +1.54%  '

# Row 10
Set-TextValue $ws.Range('E10') '  +1.62%  '

# Row 11
Set-TextValue $ws.Range('E11') '  +0.29%  '

# Row 12
Set-TextValue $ws.Range('E12') '  +2.29%  '

# Row 13
Set-TextValue $ws.Range('D13') '2.489.03'
Set-TextValue $ws.Range('E13') '  +3.09%  '

# Row 14
Set-TextValue $ws.Range('D14') '22.07'
Set-TextValue $ws.Range('E14') '  +0.12%  '

# Row 15
Set-TextValue $ws.Range('E15') '  +0.81%  '

# Row 16
Set-TextValue $ws.Range('E16') '  -0.13%  '

# Row 17
Set-TextValue $ws.Range('D17') '2.166.71'
Set-TextValue $ws.Range('E17') '  +3.10%  '

# Row 18
Set-TextValue $ws.Range('D18') '39.484.95'
Set-TextValue $ws.Range('E18') '  +1.91%  '

# Row 19
Set-TextValue $ws.Range('E19') '  +1.73%  '

# Row 20
Set-TextValue $ws.Range('D20') '72.01'
Set-TextValue $ws.Range('E20') '  +0.18%  '

# Row 21
Set-TextValue $ws.Range('E21') '  +1.18%  '

# Row 22
Set-TextValue $ws.Range('D22') '229.79'
Set-TextValue $ws.Range('E22') '  +1.00%  '

# Row 24
Set-TextValue $ws.Range('E24') '  +1.28%  '

# Row 25
Set-TextValue $ws.Range('D25') '2.29'
Set-TextValue $ws.Range('E25') '  -1.80%  '

# Row 26
Set-TextValue $ws.Range('E26') '  +1.23%  '

# Row 27
Set-TextValue $ws.Range('D27') '172.37'
Set-TextValue $ws.Range('E27') '  +0.09%  '

# Row 28
Set-TextValue $ws.Range('E28') '  +0.40%  '

# Row 29
Set-TextValue $ws.Range('D29') '19.90'
Set-TextValue $ws.Range('E29') '  +2.96%  '

# Row 30
Set-TextValue $ws.Range('E30') '  +0.21%  '

# Row 31
Set-TextValue $ws.Range('E31') '  +4.18%  '

# Row 32
Set-TextValue $ws.Range('E32') '  +2.69%  '

# Row 33
Set-TextValue $ws.Range('E33') '  +1.79%  '

# Row 34
Set-TextValue $ws.Range('E34') '  -0.55%  '

# Row 35
Set-TextValue $ws.Range('E35') '  +0.82%  '

# Row 36
Set-TextValue $ws.Range('E36') '  +0.04%  '

# Row 37
Set-TextValue $ws.Range('D37') '2.43'
Set-TextValue $ws.Range('E37') '  +1.08%  '

# Row 38
Set-TextValue $ws.Range('D38') '3.63'
Set-TextValue $ws.Range('E38') '  +1.19%  '

# Row 39
Set-TextValue $ws.Range('E39') '  -0.05%  '

# Row 40
Set-TextValue $ws.Range('D40') '102.90'
Set-TextValue $ws.Range('E40') '  +0.14%  '

# Row 41
Set-TextValue $ws.Range('E41') '  -0.79%  '

# Row 42
Set-TextValue $ws.Range('E42') '  +0.23%  '

# Row 43
Set-TextValue $ws.Range('D43') '1.526.43'
Set-TextValue $ws.Range('E43') '  +0.06%  '

# Row 44
Set-TextValue $ws.Range('E44') '  -0.08%  '

# Row 45
$ws.Range('B45').Value = 'FTXToken'
$ws.Range('C45').Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
Set-TextValue $ws.Range('D45') '4.34'
Set-TextValue $ws.Range('E45') '  +4.48%  '

# Row 46
$ws.Range('B46').Value = 'ARBITRUM'
$ws.Range('C46').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
Set-TextValue $ws.Range('D46') '1.11'
Set-TextValue $ws.Range('E46') '  +5.55%  '

# Row 47
Set-TextValue $ws.Range('D47') '7.91'
Set-TextValue $ws.Range('E47') '  +1.62%  '

# Row 48
$ws.Range('B48').Value = 'Cronos'
$ws.Range('C48').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
Set-TextValue $ws.Range('D48') '0.0927'
Set-TextValue $ws.Range('E48') '  +1.70%  '

# Row 49
$ws.Range('B49').Value = 'HuobiToken'
$ws.Range('C49').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
Set-TextValue $ws.Range('D49') '2.82'
Set-TextValue $ws.Range('E49') '  +0.31%  '

# Row 50
Set-TextValue $ws.Range('D50') '2.372.24'
Set-TextValue $ws.Range('E50') '  +3.13%  '
